# Insert a new product row ("UNICTAM 1.5 GM I.M/I.V. VIAL") into the price
# list, right after "TUSSISTOP 60 MG 20 TABS." (row 79) and before
# "VITAMIN E 400MG 24 SOFT GELATIN CAPS." (old row 80). This pushes every
# row from 80..94 down by one (81..95), bumps the running total in the
# summary row, and re-stamps the footer timestamp/page/author row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Shift the product rows 92..80 down to 93..81 (process bottom-up so
#    we never clobber a row before we've read it). Column A is a plain
#    running index (row-3) that is already correct at every row, so it
#    does not need to move. Only B (name), H (ratio text), L (qty) and
#    N (price factor) carry real data.
# ---------------------------------------------------------------------
for ($r = 92; $r -ge 80; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dst, 8).Value = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($dst, 12).Value = $ws.Cells.Item($r, 12).Value2
    $ws.Cells.Item($dst, 14).Value = $ws.Cells.Item($r, 14).Value2
}

# New row 93 (old row 92's slot) is now a normal data row -- give it the
# same B:G / H:K / L:M merges every other product row has.
$ws.Range("B93:G93").Merge()
$ws.Range("H93:K93").Merge()
$ws.Range("L93:M93").Merge()

# ---------------------------------------------------------------------
# 2) Move the totals row (old row 93, K93) down to row 94 and bump the
#    grand total by the new row's quantity (56).
# ---------------------------------------------------------------------
$oldTotal = $ws.Cells.Item(93, 11).Value2
$ws.Range("K93:N93").UnMerge()
$ws.Cells.Item(94, 11).Value = $oldTotal + 56
$ws.Range("K94:N94").Merge()

# ---------------------------------------------------------------------
# 3) Move the footer row (old row 94: timestamp / page / author) down to
#    row 95.
# ---------------------------------------------------------------------
$footerA = $ws.Cells.Item(94, 1).Value2
$footerF = $ws.Cells.Item(94, 6).Value2
$footerI = $ws.Cells.Item(94, 9).Value2

$ws.Range("A94:E94").UnMerge()
$ws.Range("F94:G94").UnMerge()
$ws.Range("I94:N94").UnMerge()
$ws.Cells.Item(94, 1).Value = $null
$ws.Cells.Item(94, 6).Value = $null
$ws.Cells.Item(94, 9).Value = $null

$ws.Cells.Item(95, 1).Value = $footerA
$ws.Cells.Item(95, 6).Value = $footerF
$ws.Cells.Item(95, 9).Value = $footerI
$ws.Range("A95:E95").Merge()
$ws.Range("F95:G95").Merge()
$ws.Range("I95:N95").Merge()

# Row heights: the totals row keeps the 25.5 height that row 93 used to
# have, and the footer's height is recalculated to 16.5 in its new spot.
$ws.Rows(94).RowHeight = 25.5
$ws.Rows(95).RowHeight = 16.5

# ---------------------------------------------------------------------
# 4) Write the new row -- UNICTAM 1.5 GM I.M/I.V. VIAL.
# ---------------------------------------------------------------------
$ws.Cells.Item(80, 2).Value = "UNICTAM 1.5 GM I.M/I.V. VIAL"
$ws.Cells.Item(80, 8).Value = "8:0"
$ws.Cells.Item(80, 12).Value = 56
$ws.Cells.Item(80, 14).Value = 1

Write-Host "Done"
